$d = $word.ActiveDocument

# 1. Title change
$d.Content.Find.Execute("Materials Engineering Project II", $true, $false, $false, $false, $false, $true, 1, $false, "Processing of Metallic Materials I", 2)

# 2. Activation date
$d.Content.Find.Execute("Ativação: 01/01/2024", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2)

# 3. Add missing space after "específicas."
$d.Content.Find.Execute("específicas.Aplicar", $true, $false, $false, $false, $false, $true, 1, $false, "específicas. Aplicar", 2)

# 5. Programa resumido paragraph
$d.Content.Find.Execute("Introdução às metodologias de projeto Design Thinking, Projeto Modelo Canvas e Lean Startup. Compreensão do uso da imaginação para a criatividade no empreendedorismo startup. Gestão da criatividade. Ciclo de vida de projeto PDCA.", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução a métodos de projeto: metodologias de projeto Design Thinking, Projeto Modelo Canvas e Lean Startup. Ciclo de vida de projeto PDCA Proposta e simulação de pequeno projeto de Engenharia. Definição do problema e formação de alternativas de solução. Estabelecimento de critérios. Escolha e avaliação de soluções. Especificação da solução. Prática de escrita científica.", 2)

# 6. Programa paragraph
$d.Content.Find.Execute("1. Metodologia de projeto focada no ser humano Design Thinking. Entendimento do duplo diamante da inovação. Etapas do Design Thinking: empatia, definição do problema, ideação, prototipação do plano e teste do produto2. Processo de melhoria contínua Kaizen. Ciclo de vida de projeto PDCA (Plan-Do-Check-Act): Planejar-Desenvolver-Checar-Agir3. Oportunidades: reconhecimento e criação de oportunidades. Uso da imaginação na criação de novos projetos4. Projeto Modelo Canvas: ideação de projeto baseado em Canvas PMC5. Aplicação do Design Thinking no projeto da inovação (produto, serviço, processo).6. Elaboração de projeto de engenharia buscando inovação, aplicando a metodologia Design Thinking. Viagem didática opcional.", $true, $false, $false, $false, $false, $true, 1, $false, "1.Introdução ao projeto em Engenharia: o que é projeto em engenharia e por que projetar? Metodologias de projeto; etapas de elaboração de projeto;2.Metodologia de projeto focada no ser humano Design Thinking. Entendimento do duplo diamante da inovação. Etapas do Design Thinking: empatia, definição do problema, ideação, prototipação do plano e teste do produto;3. Processo de melhoria contínua Kaizen. Ciclo de vida de projeto PDCA (Plan-Do-Check-Act): Planejar-Desenvolver-Checar-Agir;4.Métodos e normas para redação de textos científicos;5.Desenvolvimento de um projeto temático, compreendendo: definição do problema e formação de alternativas de solução; estabelecimento de critérios; escolha e avaliação de soluções; especificação da solução;6.Noções de planejamento e gestão de projetos; organização do tempo; técnicas para a realização de apresentações; noções de aprendizagem baseada em projetos; trabalho em grupo, equipes e times7.Tutoria de projetos", 2)

# 7. Bibliografia paragraph
$d.Content.Find.Execute("- MCCAHAN, S.; ANDERSON, P.; KORTSCHOT, M.; WEISS, P.; WOODHOUSE, K. Projetos de Engenharia: uma introdução. 1ª edição. -Rio de Janeiro: LTC, 2017.- BROCKMAN, Jay B. Introdução à Engenharia - Modelagem e solução de problemas. Rio de Janeiro: LTC, 2010.- CAVALCANTI, Carolina C.; FILATRO, Andrea C. Design Thinking na educação presencial, a distância e corporativa. São Paulo: Editora Saraiva, 2016.- IMAI, Masaaki. Gemba Kaizen: Uma abordagem de bom senso à estratégia de melhoria contínua; 2 ed. Porto Alegre: Bookman, 2014.- FINOCCHIO, José. PMC Projeto modelo Canvas, 3 ed. São Paulo: Editora Saraiva, 2020.- CAMARGO, Robson; RIBAS, Thomaz. Gestão ágil de projetos: As melhores soluções para suas necessidades. São Paulo: Editora Saraiva, 2019.- VARGAS, R. V. Manual prático do plano de projeto: Utilizando o PMBOK Guide. Rio de Janeiro: Brasport, 2018.", $true, $false, $false, $false, $false, $true, 1, $false, "- BAZZO, Walter; PEREIRA, Luiz T.V. Introdução à Engenharia, 3a. edição. Florianópolis: Editora da UFSC, 2013.- COCIAN, Luis Fernando Espinosa. Introdução à Engenharia. Porto Alegre: Bookman, 2017.- BENNETT, Ronald; MILLAM, Elaine. Liderança para engenheiros. Porto Alegre: AMGH, 2014.- ALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.- MCCAHAN, S.; ANDERSON, P.; KORTSCHOT, M.; WEISS, P.; WOODHOUSE, K. Projetos de Engenharia: uma introdução. 1ª edição. -Rio de Janeiro: LTC, 2017.", 2)

# 4. Add new docente line before "7459752 - Maria Ismenia Sodero Toledo Faria"
$target = $d.Content
$target.Find.Execute("7459752 - Maria Ismenia Sodero Toledo Faria") | Out-Null
$insertRange = $target.Duplicate
$insertRange.Collapse(1)
$insertRange.InsertBefore("3586455 - Cassius Olivio Figueiredo Terra Ruchert" + [char]11)
